$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl12"
$ws.Range("C2").Value = "Ackr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 130.955829
$ws.Range("H2").Value = 392.867487
$ws.Range("I2").Value = 0.5336535908353144
$ws.Range("J2").Value = 0.5336535908353144
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 9.688363
$ws.Range("N2").Value = 29.065089
$ws.Range("O2").Value = 0.1053077753334822
$ws.Range("P2").Value = 0.1053077753334822
$ws.Range("Q2").Value = 1268.747608317927
$ws.Range("R2").Value = 11418.72847486134
$ws.Range("S2").Value = 0.05619787244959135
$ws.Range("T2").Value = 0.05619787244959135

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl12"
$ws.Range("C3").Value = "Ackr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 130.955829
$ws.Range("H3").Value = 392.867487
$ws.Range("I3").Value = 0.5336535908353144
$ws.Range("J3").Value = 0.5336535908353144
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 71.80093133333332
$ws.Range("N3").Value = 215.402794
$ws.Range("O3").Value = 0.7804410658008428
$ws.Range("P3").Value = 0.7804410658008428
$ws.Range("Q3").Value = 9402.75048572874
$ws.Range("R3").Value = 84624.75437155867
$ws.Range("S3").Value = 0.4164851771999596
$ws.Range("T3").Value = 0.4164851771999596

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl12"
$ws.Range("C4").Value = "Ackr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 130.955829
$ws.Range("H4").Value = 392.867487
$ws.Range("I4").Value = 0.5336535908353144
$ws.Range("J4").Value = 0.5336535908353144
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.51115833333333
$ws.Range("N4").Value = 31.533475
$ws.Range("O4").Value = 0.1142511588656749
$ws.Range("P4").Value = 0.1142511588656749
$ws.Range("Q4").Value = 1376.497453291925
$ws.Range("R4").Value = 12388.47707962733
$ws.Range("S4").Value = 0.06097054118576336
$ws.Range("T4").Value = 0.06097054118576336

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl12"
$ws.Range("C5").Value = "Ackr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 66.39541
$ws.Range("H5").Value = 199.18623
$ws.Range("I5").Value = 0.2705656497465488
$ws.Range("J5").Value = 0.2705656497465488
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 9.688363
$ws.Range("N5").Value = 29.065089
$ws.Range("O5").Value = 0.1053077753334822
$ws.Range("P5").Value = 0.1053077753334822
$ws.Range("Q5").Value = 643.26283361383
$ws.Range("R5").Value = 5789.36550252447
$ws.Range("S5").Value = 0.02849266665646721
$ws.Range("T5").Value = 0.02849266665646721

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl12"
$ws.Range("C6").Value = "Ackr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 66.39541
$ws.Range("H6").Value = 199.18623
$ws.Range("I6").Value = 0.2705656497465488
$ws.Range("J6").Value = 0.2705656497465488
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 71.80093133333332
$ws.Range("N6").Value = 215.402794
$ws.Range("O6").Value = 0.7804410658008428
$ws.Range("P6").Value = 0.7804410658008428
$ws.Range("Q6").Value = 4767.252274258512
$ws.Range("R6").Value = 42905.27046832661
$ws.Range("S6").Value = 0.2111605440572941
$ws.Range("T6").Value = 0.2111605440572941

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl12"
$ws.Range("C7").Value = "Ackr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 66.39541
$ws.Range("H7").Value = 199.18623
$ws.Range("I7").Value = 0.2705656497465488
$ws.Range("J7").Value = 0.2705656497465488
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.51115833333333
$ws.Range("N7").Value = 31.533475
$ws.Range("O7").Value = 0.1142511588656749
$ws.Range("P7").Value = 0.1142511588656749
$ws.Range("Q7").Value = 697.8926671165834
$ws.Range("R7").Value = 6281.03400404925
$ws.Range("S7").Value = 0.03091243903278749
$ws.Range("T7").Value = 0.03091243903278749

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl12"
$ws.Range("C8").Value = "Ackr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 48.043585
$ws.Range("H8").Value = 144.130755
$ws.Range("I8").Value = 0.1957807594181367
$ws.Range("J8").Value = 0.1957807594181367
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 9.688363
$ws.Range("N8").Value = 29.065089
$ws.Range("O8").Value = 0.1053077753334822
$ws.Range("P8").Value = 0.1053077753334822
$ws.Range("Q8").Value = 465.463691301355
$ws.Range("R8").Value = 4189.173221712195
$ws.Range("S8").Value = 0.02061723622742368
$ws.Range("T8").Value = 0.02061723622742368

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl12"
$ws.Range("C9").Value = "Ackr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 48.043585
$ws.Range("H9").Value = 144.130755
$ws.Range("I9").Value = 0.1957807594181367
$ws.Range("J9").Value = 0.1957807594181367
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 71.80093133333332
$ws.Range("N9").Value = 215.402794
$ws.Range("O9").Value = 0.7804410658008428
$ws.Range("P9").Value = 0.7804410658008428
$ws.Range("Q9").Value = 3449.574147592163
$ws.Range("R9").Value = 31046.16732832947
$ws.Range("S9").Value = 0.152795344543589
$ws.Range("T9").Value = 0.152795344543589

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cxcl12"
$ws.Range("C10").Value = "Ackr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 48.043585
$ws.Range("H10").Value = 144.130755
$ws.Range("I10").Value = 0.1957807594181367
$ws.Range("J10").Value = 0.1957807594181367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.51115833333333
$ws.Range("N10").Value = 31.533475
$ws.Range("O10").Value = 0.1142511588656749
$ws.Range("P10").Value = 0.1142511588656749
$ws.Range("Q10").Value = 504.9937288359584
$ws.Range("R10").Value = 4544.943559523625
$ws.Range("S10").Value = 0.02236817864712401
$ws.Range("T10").Value = 0.02236817864712401
